$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" conversion summary text (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$text = $wsHoja1.Range("A1").Value()
$text = $text.Replace("1000 Bs = 4.57 = 18059.36 pesos", "1000 Bs = 4.61 = 18295.3 pesos")
$text = $text.Replace("18059.36 pesos = 4.53 = 937.62 Bs", "18295.3 pesos = 4.6 = 942.23 Bs")
$wsHoja1.Range("A1").Value = $text

# --- Update the "tasas" sheet rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 217
$wsTasas.Range("O10").Value = 3970.08
$wsTasas.Range("N12").Value = 3980.51
$wsTasas.Range("O12").Value = 205
